$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the oldest 5 years of data (2005年-2009年), which are rows 2-6.
# This shifts the remaining rows (2010年-2020年) up so they become rows 2-12.
$ws.Range("A2:J6").Delete()

# Append the newly reported year (2021年) as the new last row (row 13).
$ws.Range("A13").Value = "2021年"
$ws.Range("B13").Value = 59.2427
$ws.Range("C13").Value = 1343.2619
$ws.Range("D13").Value = 894.5577
$ws.Range("E13").Value = 216.9023
$ws.Range("F13").Value = 7.0403
$ws.Range("G13").Value = 3492.1064
$ws.Range("H13").Value = 521.9019
$ws.Range("I13").Value = 8901.9903
$ws.Range("J13").Value = 105.3572

# Match the styling used by the other year-label cells in column A
# (bold font, thin border all around, centered/top aligned) by copying
# the formatting from the row above instead of rebuilding it by hand.
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$excel.CutCopyMode = $false
